# "arrumando o grafico de 2 cores" - append the two missing data points
# (2023-11-06 and 2023-11-07) to the forerunner_245 price log so the
# 2-color chart picks them up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A40").Value = 45236
$ws.Range("A40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B40").Value = "20:40"
$ws.Range("C40").Value = 1819
$ws.Range("D40").Value = "amazon"
$ws.Range("E40").Value = "preto"

$ws.Range("A41").Value = 45237
$ws.Range("A41").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B41").Value = "20:28"
$ws.Range("C41").Value = 2026
$ws.Range("D41").Value = "amazon"
$ws.Range("E41").Value = "preto"
